# Stockwise template: swap the "Quantity" / "Unit" header labels between
# columns B and C, drop the stray empty/styled D1 cell, and leave the
# selection where the author left off (E5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the B1 / C1 header text (Количество <-> Ед. изм.) while keeping
# each cell's existing style untouched.
$b1Text = $ws.Range("B1").Value()
$c1Text = $ws.Range("C1").Value()
$ws.Range("B1").Value = $c1Text
$ws.Range("C1").Value = $b1Text

# D1 only ever carried formatting with no content - remove it outright
# so it no longer shows up as a cell record.
$ws.Range("D1").Clear()

# Reflect the author's final cursor position.
[void]$ws.Range("E5").Select()
